# Weekly data refresh: insert a new week's price record as row 28,
# shifting all subsequent rows (old 28..116) down to (29..117).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 28..116 down by one to make room for the new weekly entry.
$ws.Rows.Item(28).Insert()

# Populate the newly inserted row 28 with this week's data point.
$ws.Range("A28").Value = 11
$ws.Range("B28").Value = "Vega Monumental Concepción"
$ws.Range("C28").Value = "Bíobío"
$ws.Range("D28").Value = 45037
$ws.Range("E28").Value = 8
$ws.Range("F28").Value = 100112012
$ws.Range("G28").Value = "Espinaca"
$ws.Range("H28").Value = "Sin especificar"
$ws.Range("I28").Value = "Primera"
$ws.Range("J28").Value = 40
$ws.Range("K28").Value = 9000
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = 9500
$ws.Range("N28").Value = "$/cuna 10 kilos"
$ws.Range("O28").Value = "Región Metropolitana"
$ws.Range("P28").Value = 950
$ws.Range("Q28").Value = 10
$ws.Range("R28").Value = "Hortaliza"
